# 2025-09-30 12:38 JST scraper append: two new listings land right after the
# pinned top (highest score) row; everything else shifts down by two rows,
# and every row's "fetched at" timestamp is refreshed to the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-09-30 12:38:27"

# --- 1. Make room: insert two blank rows right after row 2 (the existing
#        rows 3-7 shift down to 5-9; row 2 - the 298-score item - stays put).
$ws.Rows("3:4").Insert()

# --- 2. Populate the two freshly inserted rows with the new listings.
$ws.Range("A3").Value = $newTimestamp
$ws.Range("B3").Value = "【急募】メモリデータ管理ツール開発のプロフェッショナル募集"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5393508"
$ws.Range("G3").Value = 158
$ws.Range("H3").Value = "◆ツール,開発 ◇管理"

$ws.Range("A4").Value = $newTimestamp
$ws.Range("B4").Value = "【RPA/Power Automate】税務システム自動化プロジェクトの依頼"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5403634"
$ws.Range("G4").Value = 103
$ws.Range("H4").Value = "◆自動化"

# --- 3. Refresh the "fetched at" timestamp on every data row (it advances
#        uniformly to the new scrape run, including rows that otherwise
#        carried forward unchanged).
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp
$ws.Range("A7").Value = $newTimestamp
$ws.Range("A8").Value = $newTimestamp
$ws.Range("A9").Value = $newTimestamp

# --- 4. Column H ("スキル概要") widened by one unit to fit the longer tags.
$ws.Columns.Item(8).ColumnWidth = 12.166666666666666

# --- 5. Rebuild the URL hyperlinks top to bottom so relationship ids stay
#        in row order (rId1..rId8), then restore the "Hyperlink" style that
#        Excel applies to linked cells on every F-column entry.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5403583")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5393508")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5403634")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5403166")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5403527")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5403384")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5403072")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5399347")

$ws.Range("F2:F9").Style = "Hyperlink"
